# Populate the "Skills" Total-Jobs figures (column F, rows 17-28) on the
# "Data" worksheet with the counts scraped for each skill. The values are
# entered as text (matching how the rest of the sheet's Total Jobs column is
# stored) and wrap text is enabled, mirroring the existing rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$skillCounts = @(8838, 5080, 1803, 3705, 253, 1246, 33, 17, 196, 104, 247, 6959)

for ($i = 0; $i -lt $skillCounts.Length; $i++) {
    $row = 17 + $i
    $cell = $ws.Cells.Item($row, 6)
    $cell.Style = "Normal"
    $cell.WrapText = $true
    $cell.Value = "'" + [string]$skillCounts[$i]
}
